$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    -5.009525469530934,
    -1.150069070561413,
    -0.06193845495403641,
    -0.4057216945391336,
    0.0133888005757794,
    0.1024167791172677,
    0.1324278584973791,
    0.02750117391858695,
    0.02290428449143454
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}
